$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'73.025.50"
$c.Style = "Normal"

$ws.Range("E2").Value = "  +2.64%  "
$c = $ws.Range("D3")
$c.Value = "'3.999.43"
$c.Style = "Normal"

$ws.Range("E3").Value = "  +1.35%  "
$c = $ws.Range("D4")
$c.Value = "'0.998"
$c.Style = "Normal"

$ws.Range("E4").Value = "  -0.24%  "
$c = $ws.Range("D5")
$c.Value = "'622.73"
$c.Style = "Normal"

$ws.Range("E5").Value = "  +16.15%  "
$c = $ws.Range("D6")
$c.Value = "'163.98"
$c.Style = "Normal"

$ws.Range("E6").Value = "  +10.66%  "
$c = $ws.Range("D7")
$c.Value = "'0.688"
$c.Style = "Normal"

$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -0.17%  "
$c = $ws.Range("D9")
$c.Value = "'0.759"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +2.32%  "
$c = $ws.Range("D10")
$c.Value = "'0.170"
$c.Style = "Normal"

$ws.Range("E10").Value = "  +2.09%  "
$c = $ws.Range("D11")
$c.Value = "'54.50"
$c.Style = "Normal"

$ws.Range("E11").Value = "  -0.51%  "
$c = $ws.Range("D12")
$c.Value = "'0.0000322"
$c.Style = "Normal"

$ws.Range("E12").Value = "  +0.83%  "
$c = $ws.Range("D13")
$c.Value = "'11.16"
$c.Style = "Normal"

$ws.Range("E13").Value = "  +4.58%  "
$c = $ws.Range("D14")
$c.Value = "'4.618.84"
$c.Style = "Normal"

$ws.Range("E14").Value = "  +0.87%  "
$c = $ws.Range("D15")
$c.Value = "'3.980.39"
$c.Style = "Normal"

$ws.Range("E15").Value = "  +0.88%  "
$c = $ws.Range("D16")
$c.Value = "'1.27"
$c.Style = "Normal"

$ws.Range("E16").Value = "  +8.96%  "
$c = $ws.Range("D17")
$c.Value = "'14.21"
$c.Style = "Normal"

$ws.Range("E17").Value = "  +1.61%  "
$c = $ws.Range("D18")
$c.Value = "'20.80"
$c.Style = "Normal"

$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +0.46%  "
$c = $ws.Range("D20")
$c.Value = "'72.599.92"
$c.Style = "Normal"

$ws.Range("E20").Value = "  +2.17%  "
$c = $ws.Range("D21")
$c.Value = "'442.34"
$c.Style = "Normal"

$ws.Range("E21").Value = "  +3.41%  "
$c = $ws.Range("D22")
$c.Value = "'4.91"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +16.05%  "
$c = $ws.Range("D23")
$c.Value = "'96.98"
$c.Style = "Normal"

$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  -2.97%  "
$c = $ws.Range("D25")
$c.Value = "'14.52"
$c.Style = "Normal"

$ws.Range("E25").Value = "  -0.40%  "
$c = $ws.Range("D26")
$c.Value = "'4.32"
$c.Style = "Normal"

$ws.Range("E26").Value = "  +10.78%  "
$c = $ws.Range("D27")
$c.Value = "'11.46"
$c.Style = "Normal"

$ws.Range("E27").Value = "  +1.34%  "
$c = $ws.Range("B28")
$c.Value = "'Filecoin"
$c.Style = "Normal"

$c = $ws.Range("C28")
$c.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'10.60"
$c.Style = "Normal"

$ws.Range("E28").Value = "  -1.58%  "
$c = $ws.Range("B29")
$c.Value = "'LEO"
$c.Style = "Normal"

$c = $ws.Range("C29")
$c.Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'5.97"
$c.Style = "Normal"

$ws.Range("E29").Value = "  +1.00%  "
$c = $ws.Range("D30")
$c.Value = "'36.63"
$c.Style = "Normal"

$ws.Range("E30").Value = "  +0.21%  "
$c = $ws.Range("D31")
$c.Value = "'7.78"
$c.Style = "Normal"

$ws.Range("E31").Value = "  -0.80%  "
$c = $ws.Range("D32")
$c.Value = "'14.00"
$c.Style = "Normal"

$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("E33").Value = "  +0.42%  "
$c = $ws.Range("D34")
$c.Value = "'72.65"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +10.81%  "
$c = $ws.Range("D35")
$c.Value = "'48.35"
$c.Style = "Normal"

$ws.Range("E35").Value = "  -4.78%  "
$c = $ws.Range("D36")
$c.Value = "'657.97"
$c.Style = "Normal"

$ws.Range("E36").Value = "  -3.37%  "
$c = $ws.Range("D37")
$c.Value = "'0.0₃0919"
$c.Style = "Normal"

$ws.Range("E37").Value = "  +12.95%  "
$c = $ws.Range("D38")
$c.Value = "'0.442"
$c.Style = "Normal"

$ws.Range("E38").Value = "  +0.11%  "
$c = $ws.Range("B39")
$c.Value = "'Kaspa"
$c.Style = "Normal"

$c = $ws.Range("C39")
$c.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.148"
$c.Style = "Normal"

$ws.Range("E39").Value = "  -0.57%  "
$c = $ws.Range("B40")
$c.Value = "'ThetaToken"
$c.Style = "Normal"

$c = $ws.Range("C40")
$c.Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'3.39"
$c.Style = "Normal"

$ws.Range("E40").Value = "  +0.45%  "
$c = $ws.Range("D41")
$c.Value = "'0.999"
$c.Style = "Normal"

$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("E43").Value = "  +0.19%  "
$c = $ws.Range("D44")
$c.Value = "'0.0492"
$c.Style = "Normal"

$ws.Range("E44").Value = "  +1.85%  "
$c = $ws.Range("D45")
$c.Value = "'10.68"
$c.Style = "Normal"

$ws.Range("E45").Value = "  +4.05%  "
$c = $ws.Range("D46")
$c.Value = "'0.151"
$c.Style = "Normal"

$ws.Range("E46").Value = "  +1.47%  "
$c = $ws.Range("D47")
$c.Value = "'2.67"
$c.Style = "Normal"

$ws.Range("E47").Value = "  +0.73%  "
$c = $ws.Range("D48")
$c.Value = "'3.42"
$c.Style = "Normal"

$ws.Range("E48").Value = "  +1.97%  "
$c = $ws.Range("B49")
$c.Value = "'Stacks"
$c.Style = "Normal"

$c = $ws.Range("C49")
$c.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'3.10"
$c.Style = "Normal"

$ws.Range("E49").Value = "  +3.39%  "
$c = $ws.Range("B50")
$c.Value = "'Maker"
$c.Style = "Normal"

$c = $ws.Range("C50")
$c.Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'2.908.27"
$c.Style = "Normal"

$ws.Range("E50").Value = "  +11.68%  "
$c = $ws.Range("D51")
$c.Value = "'3.42"
$c.Style = "Normal"

$ws.Range("E51").Value = "  +5.15%  "